# Assignment 1.2 - add the GitHub repo link paragraph and flag the
# (now-pasted) screenshot's run as "no proofing" the same way Word does
# once a picture has been dropped into the body.
#
# Corresponds to the commit:
#   "Added github link to word documents for assignments 1.2 and 1.3"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) The blank paragraph right after "Assignment 1.2" gets filled in with
#    the GitHub link text.
# ---------------------------------------------------------------------
$anchorText = $d.Content
$found = $anchorText.Find.Execute("Assignment 1.2")

if ($found) {
    # Walk the paragraph collection to find the paragraph that contains the
    # match, then target the paragraph immediately following it (this is
    # more robust than assuming a fixed paragraph index).
    $targetIndex = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs.Item($i)
        if (($anchorText.Start -ge $para.Range.Start) -and ($anchorText.Start -lt $para.Range.End)) {
            $targetIndex = $i + 1
        }
    }

    if ($targetIndex -ge 1 -and $targetIndex -le $d.Paragraphs.Count) {
        $targetPara = $d.Paragraphs.Item($targetIndex)

        # A paragraph's Range.Text always includes the trailing paragraph
        # mark (\r), so trim it off before checking for emptiness. Only
        # fill it in if it's still empty (idempotent / safe to re-run).
        $existingText = $targetPara.Range.Text
        if ($existingText -eq $null) { $existingText = "" }
        if ($existingText.Trim([char]13, [char]10) -eq "") {
            $insertionPoint = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)
            $insertionPoint.InsertAfter("Github link: https://github.com/Sara-Renee/csd-340/")
        }
    }
}

# ---------------------------------------------------------------------
# 2) The second inline picture (the full-page code screenshot) picks up
#    the same "noProof" run formatting the first picture already has.
# ---------------------------------------------------------------------
if ($d.InlineShapes.Count -ge 2) {
    $secondPicture = $d.InlineShapes.Item(2)
    $secondPicture.Range.NoProofing = 1
}
